# Created experiment order generation script
# Re-running the generator reshuffles task-order rows for each condition
# sheet and re-randomizes the stimulus-order tab names/order.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: rewrite each sheet's data table (rows 2..N of columns A:B).
# Addressed by ORIGINAL physical position (1..5) before any tab
# reordering happens, since Worksheets.Item(<index>) tracks tab order.
# ---------------------------------------------------------------------

function Set-TaskOrderRows {
    param(
        $ws,
        [string[]]$values
    )

    $oldLastRow = $ws.UsedRange.Rows.Count
    if ($oldLastRow -lt 1) { $oldLastRow = 1 }
    $newLastRow = $values.Count + 1

    # Grow: paste the bold/bordered/centered "A" style from A2 down to
    # any newly needed rows before filling in values.
    if ($newLastRow -gt $oldLastRow) {
        $ws.Range("A2").Copy()
        $ws.Range($ws.Cells.Item($oldLastRow + 1, 1), $ws.Cells.Item($newLastRow, 1)).PasteSpecial(-4122)
    }

    for ($i = 0; $i -lt $values.Count; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 1).Value = $i
        $ws.Cells.Item($row, 2).Value = $values[$i]
    }

    # Shrink: drop any now-unused trailing rows entirely.
    if ($newLastRow -lt $oldLastRow) {
        $ws.Range($ws.Cells.Item($newLastRow + 1, 1), $ws.Cells.Item($oldLastRow, 2)).EntireRow.Delete()
    }
}

# sheet1.xml (tab 1 originally "GNG_TO-...") -> becomes the "NB" task order
$wsA = $wb.Worksheets.Item(1)
Set-TaskOrderRows $wsA @(
    "TB-16515890217845197.csv",
    "OB-16515890206829958.csv",
    "OB-16515890205557156.csv",
    "TB-16515890219585898.csv",
    "ZB-match_8-1651589020140042.csv",
    "ZB-match_0-1651589019852792.csv",
    "TB-16515890213916445.csv",
    "OB-16515890207475772.csv",
    "ZB-match_6-16515890200932014.csv"
)

# sheet2.xml (tab 2 originally "NB_TO-...") -> becomes the "TOL" task order
$wsB = $wb.Worksheets.Item(2)
Set-TaskOrderRows $wsB @(
    "MM_stims-16515890219898074.csv",
    "ZM_stims-16515890219742167.csv",
    "MM_stims-16515890220054326.csv",
    "ZM_stims-16515890219898074.csv",
    "MM_stims-16515890220210583.csv",
    "ZM_stims-16515890220054326.csv"
)

# sheet3.xml (tab 3 "RS_TO-...") -> resting-state task order is unchanged,
# only its generated tab name changes (handled in the rename step below).

# sheet4.xml (tab 4 originally "TOL_TO-...") -> becomes the "GNG" task order
$wsD = $wb.Worksheets.Item(4)
Set-TaskOrderRows $wsD @(
    "go_stims-16515890220210583.csv",
    "GNG_stims-16515890220366833.csv",
    "go_stims-16515890220366833.csv",
    "GNG_stims-16515890220523074.csv"
)

# sheet5.xml (tab 5 "vSAT_TO-...") -> vSAT task order refreshed, same length
$wsE = $wb.Worksheets.Item(5)
Set-TaskOrderRows $wsE @(
    "vSAT_stims-16515890220835586.csv",
    "vSAT_stims-16515890220991824.csv",
    "SAT_stims-16515890220523074.csv",
    "vSAT_stims-16515890220679328.csv"
)

# ---------------------------------------------------------------------
# Step 2: rename every tab to the newly generated run id.
# ---------------------------------------------------------------------

$wsA.Name = "NB_TO-16515890219742167"
$wsB.Name = "TOL_TO-16515890220210583"
$wb.Worksheets.Item(3).Name = "RS_TO-16515890220210583"
$wsD.Name = "GNG_TO-16515890220523074"
$wsE.Name = "vSAT_TO-1651589022114808"

# ---------------------------------------------------------------------
# Step 3: reorder the tabs to NB, TOL, RS, GNG, vSAT by repeatedly
# moving the next desired tab to the end of the tab strip.
# ---------------------------------------------------------------------

$orderedNames = @(
    "NB_TO-16515890219742167",
    "TOL_TO-16515890220210583",
    "RS_TO-16515890220210583",
    "GNG_TO-16515890220523074",
    "vSAT_TO-1651589022114808"
)

foreach ($name in $orderedNames) {
    $sheet = $wb.Worksheets.Item($name)
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    if ($sheet.Name -ne $lastSheet.Name) {
        $sheet.Move($null, $lastSheet)
    }
}
